$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target column is CS = column 97 (one past existing CR = 96)
$col = 97

# Set the new column width to match the existing data columns (raw XML width = 12)
$ws.Columns.Item($col).ColumnWidth = 11.17

# --- Row 1: header date, stored as literal text (not an auto-converted date) ---
# A plain Value="2024/12/14" gets auto-parsed into a date serial by the COM layer
# (and mints a brand-new number-format style). Instead, compute the text via a
# formula in a scratch cell (formula results are never re-parsed as dates), paste
# only the resulting value into the header cell, then paste the CR1 cell's
# formatting on top so it reuses the existing style (s="1") instead of creating one.
$scratch = $ws.Cells.Item(60, $col)
$scratch.Formula = '="2024/12/14"'
$h = $ws.Cells.Item(1, $col)
$scratch.Copy()
$h.PasteSpecial(-4163)          # xlPasteValues
$scratch.Clear()

$hsrc = $ws.Cells.Item(1, 96)   # CR1, style 1 (plain header style)
$hsrc.Copy()
$h.PasteSpecial(-4122)          # xlPasteFormats: reuse CR1's existing style

# --- Data rows 2-53 ---
# For each row, a same-style source cell already exists somewhere in column CR (96);
# picking one per style lets PasteSpecial(xlPasteFormats) reuse the existing cellXfs entry
# instead of minting a brand-new style for each write.
$styleSource = @{ 1 = 1; 2 = 47; 3 = 2 }   # style -> a CR-column row known to carry that style

$data = @(
    @{Row=2; Value=125; Style=3}
    @{Row=3; Value=137.7; Style=3}
    @{Row=4; Value=173.5; Style=1}
    @{Row=5; Value=163.4; Style=1}
    @{Row=6; Value=135.9; Style=3}
    @{Row=7; Value=204.2; Style=1}
    @{Row=8; Value=157.1; Style=1}
    @{Row=9; Value=171.4; Style=1}
    @{Row=10; Value=165.6; Style=1}
    @{Row=11; Value=184.1; Style=1}
    @{Row=12; Value=164.9; Style=1}
    @{Row=13; Value=147.5; Style=1}
    @{Row=14; Value=116.5; Style=2}
    @{Row=15; Value=133.3; Style=3}
    @{Row=16; Value=158.5; Style=1}
    @{Row=17; Value=116.4; Style=2}
    @{Row=18; Value=124.5; Style=2}
    @{Row=19; Value=220.1; Style=1}
    @{Row=20; Value=153.2; Style=1}
    @{Row=21; Value=139.4; Style=3}
    @{Row=22; Value=177; Style=1}
    @{Row=23; Value=277.6; Style=1}
    @{Row=24; Value=115.2; Style=2}
    @{Row=25; Value=145.2; Style=1}
    @{Row=26; Value=134.1; Style=3}
    @{Row=27; Value=202.6; Style=1}
    @{Row=28; Value=159; Style=1}
    @{Row=29; Value=149.9; Style=1}
    @{Row=30; Value=150.1; Style=1}
    @{Row=31; Value=181.3; Style=1}
    @{Row=32; Value=194; Style=1}
    @{Row=33; Value=138.7; Style=3}
    @{Row=34; Value=221.6; Style=1}
    @{Row=35; Value=261.4; Style=1}
    @{Row=36; Value=154.3; Style=1}
    @{Row=37; Value=151.7; Style=1}
    @{Row=38; Value=179.6; Style=1}
    @{Row=39; Value=161.2; Style=1}
    @{Row=40; Value=142.9; Style=1}
    @{Row=41; Value=149.4; Style=1}
    @{Row=42; Value=145.4; Style=1}
    @{Row=43; Value=152.4; Style=1}
    @{Row=44; Value=133.4; Style=3}
    @{Row=45; Value=211.4; Style=1}
    @{Row=46; Value=150; Style=1}
    @{Row=47; Value=188.7; Style=1}
    @{Row=48; Value=137.6; Style=3}
    @{Row=49; Value=110.8; Style=2}
    @{Row=50; Value=234.7; Style=1}
    @{Row=51; Value=181.5; Style=1}
    @{Row=52; Value=167.3; Style=1}
    @{Row=53; Value=186.3; Style=1}
)

foreach ($d in $data) {
    $src = $ws.Cells.Item($styleSource[$d.Style], 96)
    $dst = $ws.Cells.Item($d.Row, $col)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $d.Value
}

$excel.CutCopyMode = 0
